$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9422723054885864
$ws.Range("B1").Value = 1.697656154632568
$ws.Range("C1").Value = 4.580130577087402
$ws.Range("D1").Value = 1.896753549575806
$ws.Range("E1").Value = 1.037350416183472
